# Sprint 2 Report.pptx edit
#
# Adds a new bullet to the "Changes For Next Sprint" slide (slide 4):
#   "Convention of testing, and naming conventions should be consistently
#    stated before sprint begins."
# inserted as a new paragraph right after the existing "Testing integrated
# more within development cycle ..." bullet, and before the trailing blank
# paragraph.

$p = $ppt.ActivePresentation

# Find the "Changes For Next Sprint" slide instead of hard-coding its index.
$targetSlide = $null
for ($i = 1; $i -le $p.Slides.Count; $i++) {
    $candidate = $p.Slides.Item($i)
    if ($candidate.Shapes.Count -ge 1 -and $candidate.Shapes.Item(1).HasTextFrame) {
        if ($candidate.Shapes.Item(1).TextFrame.TextRange.Text -eq "Changes For Next Sprint") {
            $targetSlide = $candidate
            break
        }
    }
}

if ($targetSlide -eq $null) {
    $targetSlide = $p.Slides.Item($p.Slides.Count)
}

# Find the body placeholder that holds the bullet list (the shape that is
# not the title).
$body = $null
for ($j = 1; $j -le $targetSlide.Shapes.Count; $j++) {
    $candidateShape = $targetSlide.Shapes.Item($j)
    if ($candidateShape.HasTextFrame -and $candidateShape.TextFrame.TextRange.Text -ne "Changes For Next Sprint") {
        $body = $candidateShape
        break
    }
}

$tr = $body.TextFrame.TextRange

# Locate the paragraph that starts with "Testing integrated" so the new
# bullet is inserted right after it (and therefore right before the
# trailing empty paragraph). NOTE: `.Index` is not reliable on the
# sub-ranges returned by `Paragraphs()` in this host, so the 1-based
# position is tracked manually via the loop counter instead.
$paraCount = $tr.Paragraphs().Count
$anchorIdx = -1
for ($k = 1; $k -le $paraCount; $k++) {
    $candidatePara = $tr.Paragraphs($k)
    if ($candidatePara.Text -like "Testing integrated*") {
        $anchorIdx = $k
    }
}

if ($anchorIdx -eq -1) {
    # Fall back to the last non-empty paragraph.
    for ($k = 1; $k -le $paraCount; $k++) {
        $candidatePara = $tr.Paragraphs($k)
        if ($candidatePara.Text.Length -gt 0) {
            $anchorIdx = $k
        }
    }
}

$anchorPara = $tr.Paragraphs($anchorIdx)

$firstRunText = "Convention of testing, and naming conventions should be consistently stated "
$secondRunText = "before sprint begins."

# Insert a new paragraph (leading carriage return) right after the anchor
# paragraph, containing both sentences.
$null = $anchorPara.InsertAfter("`r" + $firstRunText + $secondRunText)

# Re-fetch the freshly created paragraph: it now immediately follows the
# anchor paragraph, i.e. at position (anchorIdx + 1).
$newParaIndex = $anchorIdx + 1
$newPara = $tr.Paragraphs($newParaIndex)

# Make sure the new bullet uses the same 16pt size as its neighbours.
$newPara.Font.Size = 16

# Split the new paragraph into two runs (matching the source document,
# where the two sentences ended up as separate runs) by nudging the
# formatting of the second sentence only.
$secondRunStart = $newPara.Start + $firstRunText.Length
$secondRun = $tr.Characters($secondRunStart, $secondRunText.Length)
$secondRun.Font.Size = 16
